# Update the "Marking" (per-correct-answer weight) and resulting "Total"
# score on the quiz marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking weight for a correct answer: 3 -> 5
$ws.Range("B11").Value = 5

# Total correct-answer marks recalculated with the new weight: 48 -> 80
$ws.Range("B12").Value = 80

# Displayed "scored/total" fraction: 39/84 -> 80/140
$ws.Range("E12").Value = "80/140"
